$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The two observation records (row 2 and row 4) had their Id (A), Antal (I),
# Ost (Q) and Nord (R) values swapped between each other.

# Capture original values (use Value2 - Value is unreliable in this runtime)
$a2 = $ws.Range("A2").Value2
$i2 = $ws.Range("I2").Value2
$q2 = $ws.Range("Q2").Value2
$r2 = $ws.Range("R2").Value2

$a4 = $ws.Range("A4").Value2
$i4 = $ws.Range("I4").Value2
$q4 = $ws.Range("Q4").Value2
$r4 = $ws.Range("R4").Value2

# Antal ("I") column is stored as text in this workbook, not a number.
# Force text formatting so the swapped numeric-looking strings stay text.
$ws.Range("I2").NumberFormat = "@"
$ws.Range("I4").NumberFormat = "@"

# Row 2 gets row 4's original values
$ws.Range("A2").Value2 = $a4
$ws.Range("I2").Value2 = "$i4"
$ws.Range("Q2").Value2 = $q4
$ws.Range("R2").Value2 = $r4

# Row 4 gets row 2's original values
$ws.Range("A4").Value2 = $a2
$ws.Range("I4").Value2 = "$i2"
$ws.Range("Q4").Value2 = $q2
$ws.Range("R4").Value2 = $r2
